# Generate Report for Handback
#
# This mirrors a localization "handback" report generation run: the
# latest target + handback xliff for "a.md" are now available for both
# the zh-cn and de-de locales, so the Overview / per-locale sheets are
# refreshed with the new status text, the new handback file links and
# the new handback timestamps.

$wb = $excel.ActiveWorkbook

$statusText   = "Handed back: in sync with en-US"
$aMdUrl       = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4dd47b1923181cbd050d11c94e46ba5693cfa402/e2e/a.md"
$zhHandback   = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deHandback   = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$zhHandbackDt = "2016-08-28 10:44:32"
$deHandbackDt = "2016-08-28 10:44:39"

# ---------------------------------------------------------------------
# Overview sheet: refresh the status column for both locales/rows
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.144371396019366
$wsOverview.Columns.Item(6).ColumnWidth = 29.144371396019366

# ---------------------------------------------------------------------
# zh-cn sheet: a.md has been handed back
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("I2").Value = "a.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $aMdUrl, "", "", "a.md")
$wsZh.Range("J2").Value = $zhHandback
$wsZh.Range("K2").Value = $zhHandbackDt

$wsZh.Range("I3").Value = "a.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $aMdUrl, "", "", "a.md")
$wsZh.Range("J3").Value = $zhHandback
$wsZh.Range("K3").Value = $zhHandbackDt

$wsZh.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet: a.md has been handed back
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("I2").Value = "a.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $aMdUrl, "", "", "a.md")
$wsDe.Range("J2").Value = $deHandback
$wsDe.Range("K2").Value = $deHandbackDt

$wsDe.Range("I3").Value = "a.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $aMdUrl, "", "", "a.md")
$wsDe.Range("J3").Value = $deHandback
$wsDe.Range("K3").Value = $deHandbackDt

$wsDe.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664
